$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-09 Monday", 2) | Out-Null
$d.Content.Find.Execute("756×9=6804", $true, $false, $false, $false, $false, $true, 1, $false, "790×2=1580", 2) | Out-Null
$d.Content.Find.Execute("928×9=8352", $true, $false, $false, $false, $false, $true, 1, $false, "997×2=1994", 2) | Out-Null
$d.Content.Find.Execute("358×5=1790", $true, $false, $false, $false, $false, $true, 1, $false, "120×5=600", 2) | Out-Null
$d.Content.Find.Execute("681×4=2724", $true, $false, $false, $false, $false, $true, 1, $false, "728×3=2184", 2) | Out-Null
$d.Content.Find.Execute("431×9=3879", $true, $false, $false, $false, $false, $true, 1, $false, "699×4=2796", 2) | Out-Null
$d.Content.Find.Execute("253×3=759", $true, $false, $false, $false, $false, $true, 1, $false, "430×6=2580", 2) | Out-Null
$d.Content.Find.Execute("904×7=6328", $true, $false, $false, $false, $false, $true, 1, $false, "738×4=2952", 2) | Out-Null
$d.Content.Find.Execute("801×7=5607", $true, $false, $false, $false, $false, $true, 1, $false, "532×2=1064", 2) | Out-Null
$d.Content.Find.Execute("447×5=2235", $true, $false, $false, $false, $false, $true, 1, $false, "711×5=3555", 2) | Out-Null
$d.Content.Find.Execute("920×2=1840", $true, $false, $false, $false, $false, $true, 1, $false, "659×6=3954", 2) | Out-Null
$d.Content.Find.Execute("322×4=1288", $true, $false, $false, $false, $false, $true, 1, $false, "343×8=2744", 2) | Out-Null
$d.Content.Find.Execute("118×2=236", $true, $false, $false, $false, $false, $true, 1, $false, "728×6=4368", 2) | Out-Null
$d.Content.Find.Execute("910×9=8190", $true, $false, $false, $false, $false, $true, 1, $false, "405×4=1620", 2) | Out-Null
$d.Content.Find.Execute("948×3=2844", $true, $false, $false, $false, $false, $true, 1, $false, "641×5=3205", 2) | Out-Null
$d.Content.Find.Execute("922×4=3688", $true, $false, $false, $false, $false, $true, 1, $false, "611×3=1833", 2) | Out-Null
$d.Content.Find.Execute("740×2=1480", $true, $false, $false, $false, $false, $true, 1, $false, "927×6=5562", 2) | Out-Null
$d.Content.Find.Execute("255×3=765", $true, $false, $false, $false, $false, $true, 1, $false, "382×2=764", 2) | Out-Null
$d.Content.Find.Execute("874×9=7866", $true, $false, $false, $false, $false, $true, 1, $false, "350×3=1050", 2) | Out-Null
$d.Content.Find.Execute("298×4=1192", $true, $false, $false, $false, $false, $true, 1, $false, "992×7=6944", 2) | Out-Null
$d.Content.Find.Execute("169×6=1014", $true, $false, $false, $false, $false, $true, 1, $false, "798×5=3990", 2) | Out-Null
$d.Content.Find.Execute("642×7=4494", $true, $false, $false, $false, $false, $true, 1, $false, "147×6=882", 2) | Out-Null
$d.Content.Find.Execute("356×3=1068", $true, $false, $false, $false, $false, $true, 1, $false, "243×6=1458", 2) | Out-Null
$d.Content.Find.Execute("190×6=1140", $true, $false, $false, $false, $false, $true, 1, $false, "664×9=5976", 2) | Out-Null
$d.Content.Find.Execute("485×9=4365", $true, $false, $false, $false, $false, $true, 1, $false, "163×7=1141", 2) | Out-Null
$d.Content.Find.Execute("450×7=3150", $true, $false, $false, $false, $false, $true, 1, $false, "709×3=2127", 2) | Out-Null
